$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 372 (shifts existing 372:406 down to 374:408)
$ws.Rows("372:373").Insert()

# New row 372 - "Primera" quality entry for the newest reporting date
$ws.Cells.Item(372,1).Value = 7
$ws.Cells.Item(372,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(372,3).Value = "Ñuble"
$ws.Cells.Item(372,4).Value = 45132
$ws.Cells.Item(372,5).Value = 16
$ws.Cells.Item(372,6).Value = 100112017
$ws.Cells.Item(372,7).Value = "Apio"
$ws.Cells.Item(372,8).Value = "Americana (o)"
$ws.Cells.Item(372,9).Value = "Primera"
$ws.Cells.Item(372,10).Value = 180
$ws.Cells.Item(372,11).Value = 7000
$ws.Cells.Item(372,12).Value = 7000
$ws.Cells.Item(372,13).Value = 7000
$ws.Cells.Item(372,14).Value = "$/docena de matas"
$ws.Cells.Item(372,15).Value = "Provincia del Elquí"
$ws.Cells.Item(372,16).Value = 1167
$ws.Cells.Item(372,17).Value = 6
$ws.Cells.Item(372,18).Value = "Hortaliza"

# New row 373 - "Segunda" quality entry for the newest reporting date
$ws.Cells.Item(373,1).Value = 7
$ws.Cells.Item(373,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(373,3).Value = "Ñuble"
$ws.Cells.Item(373,4).Value = 45132
$ws.Cells.Item(373,5).Value = 16
$ws.Cells.Item(373,6).Value = 100112017
$ws.Cells.Item(373,7).Value = "Apio"
$ws.Cells.Item(373,8).Value = "Americana (o)"
$ws.Cells.Item(373,9).Value = "Segunda"
$ws.Cells.Item(373,10).Value = 180
$ws.Cells.Item(373,11).Value = 6000
$ws.Cells.Item(373,12).Value = 6000
$ws.Cells.Item(373,13).Value = 6000
$ws.Cells.Item(373,14).Value = "$/docena de matas"
$ws.Cells.Item(373,15).Value = "Provincia del Elquí"
$ws.Cells.Item(373,16).Value = 1000
$ws.Cells.Item(373,17).Value = 6
$ws.Cells.Item(373,18).Value = "Hortaliza"
